$d = $word.ActiveDocument

$replacements = @(
    @{old="353×2="; new="675×5="},
    @{old="955×2="; new="414×5="},
    @{old="430×4="; new="668×7="},
    @{old="634×8="; new="744×8="},
    @{old="975×4="; new="799×6="},
    @{old="987×4="; new="210×2="},
    @{old="324×4="; new="350×9="},
    @{old="296×7="; new="936×6="},
    @{old="592×2="; new="156×6="},
    @{old="265×2="; new="659×3="},
    @{old="255×8="; new="493×4="},
    @{old="446×4="; new="706×6="},
    @{old="916×9="; new="231×8="},
    @{old="268×9="; new="236×4="},
    @{old="488×3="; new="279×3="},
    @{old="524×4="; new="339×2="},
    @{old="907×3="; new="716×4="},
    @{old="688×5="; new="877×3="},
    @{old="336×2="; new="630×7="},
    @{old="243×8="; new="807×3="},
    @{old="492×4="; new="747×9="},
    @{old="648×5="; new="333×2="},
    @{old="410×7="; new="760×2="},
    @{old="699×8="; new="377×9="},
    @{old="981×4="; new="593×5="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
